# stop btn, fix isAdmin, write user info
# Append the new expense rows (17.11.2022 - 19.11.2022) below the existing
# data table on the active sheet, rows 6..16, columns A..G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(6,  "17.11.2022", "45",     "Общее",   "",                "Топливо",             "Novisa",      ""),
    @(7,  "17.11.2022", "34",     "Бригада", "",                "Жилье",               "Karpacz",     "Бригада Лелюка"),
    @(8,  "17.11.2022", "54",     "Бригада", "",                "Топливо",             "Skysawa",     "Бригада Миши"),
    @(9,  "17.11.2022", "56",     "Люди",    "Lelyuk Alexandr", "Материал",            "MCM project", ""),
    @(10, "17.11.2022", "53,89",  "Общее",   "",                "Топливо",             "Karpacz",     ""),
    @(11, "17.11.2022", "32,78",  "Общее",   "",                "Проезд - билеты",     "Karpacz",     ""),
    @(12, "17.11.2022", "43,98",  "Общее",   "",                "Инструмент",          "MCM project", ""),
    @(13, "17.11.2022", "23,69",  "Общее",   "",                "Зарплата",            "Office",      ""),
    @(14, "17.11.2022", "214,36", "Общее",   "",                "коллекция покемонов", "Office",      ""),
    @(15, "17.11.2022", "54",     "Бригада", "",                "Габилен",             "Karpacz",     "Бригада Игоря"),
    @(16, "19.11.2022", "43,89",  "Люди",    "Владислав",       "Жилье",               "Karpacz",     "")
)

# Force the newly written range to be stored as plain text (same as the
# existing rows) so date-looking and number-looking strings ("17.11.2022",
# "45", ...) are not reinterpreted as real dates/numbers, then restore the
# General format so the new cells share the same style as the rest of the
# table.
$target = $ws.Range("A6:G16")
$target.NumberFormat = "@"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    if ($row[4] -ne "") {
        $ws.Cells.Item($r, 4).Value = $row[4]
    }
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    if ($row[7] -ne "") {
        $ws.Cells.Item($r, 7).Value = $row[7]
    }
}

$target.NumberFormat = "General"
